$wb = $excel.ActiveWorkbook
$entities = $wb.Worksheets.Item("entities")
$attributes = $wb.Worksheets.Item("attributes")

# --- Write new cell text in the exact order needed so the shared-string
# table grows with the same index assignment as the target workbook. ---

# 1 (idx143): attributes!A45 = xcomputedint
$attributes.Range("A45").Value = "xcomputedint"

# 2 (idx144): attributes!P1 = expression (new column header)
$attributes.Range("P1").Value = "expression"

# 3 (idx145): attributes!A44 = xcomputedxref
$attributes.Range("A44").Value = "xcomputedxref"

# 4 (idx146): entities!D1 = abstract (new column header)
$entities.Range("D1").Value = "abstract"

# 5 (idx147): entities!A4 = Location
$entities.Range("A4").Value = "Location"

# reuse idx147 "Location"
$attributes.Range("D44").Value = "Location"
$attributes.Range("B46").Value = "Location"
$attributes.Range("B47").Value = "Location"

# 6 (idx148): attributes!P44 = {Chromosome: xstring, Position: xint}
$attributes.Range("P44").Value = "{Chromosome: xstring, Position: xint}"

# 7 (idx149): entities!C4 = entity for the computed attributes
$entities.Range("C4").Value = "entity for the computed attributes"

# 8 (idx150): attributes!A46 = Chromosome
$attributes.Range("A46").Value = "Chromosome"

# 9 (idx151): attributes!A47 = Position
$attributes.Range("A47").Value = "Position"

# --- Remaining cells in new rows that reuse existing shared strings ---

# row44 : xcomputedxref
$attributes.Range("B44").Value = "TypeTest"
$attributes.Range("C44").Value = "xref"
$attributes.Range("E44").Value = $false
$attributes.Range("F44").Value = $true
$attributes.Range("L44").Value = $true

# row45 : xcomputedint
$attributes.Range("B45").Value = "TypeTest"
$attributes.Range("C45").Value = "int"
$attributes.Range("E45").Value = $false
$attributes.Range("F45").Value = $true
$attributes.Range("L45").Value = $true
$attributes.Range("P45").Value = "xint"

# row46 : Chromosome
$attributes.Range("C46").Value = "string"
$attributes.Range("E46").Value = $false
$attributes.Range("F46").Value = $false

# row47 : Position
$attributes.Range("C47").Value = "int"
$attributes.Range("E47").Value = $true
$attributes.Range("F47").Value = $false

# --- Page setup (entities sheet gains a pageSetup element matching attributes') ---
$entities.PageSetup.PaperSize = 9

# --- Column width adjustments ---
$entities.Columns.Item(1).ColumnWidth = 14.998697916666666
$entities.Columns.Item(3).ColumnWidth = 26.166666666666668
$attributes.Columns.Item(2).ColumnWidth = 14.998697916666666
$attributes.Columns.Item(16).ColumnWidth = 28.166666666666668

# --- Selections / active sheet / view state ---
$typetest = $wb.Worksheets.Item("TypeTest")
$typetest.Activate()
$typetest.Range("D36").Select()

$entities.Activate()
$entities.Range("B5").Select()

$attributes.Activate()
$attributes.Range("F36").Select()

Write-Output "done"
